$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Append the 9 new machine rows (22-30) - matches rows "Machine 21".."Machine 29"
#    Columns left as-is for now (A,B,D,E,F,G,H,I,J,K); column C (mac_address)
#    is populated afterwards in a separate pass together with the existing
#    rows' mac addresses (mirrors how the source workbook's shared-string
#    table ends up ordered: new name/serial/ip strings first, then every mac
#    address - old format - replaced by new dash-formatted ones, in row order).
# ---------------------------------------------------------------------------
$newRows = @(
    @{r=22; id=10021; name="Machine 21"; serial="FB5962911653"; ip="192.168.0.874"},
    @{r=23; id=10022; name="Machine 22"; serial="FB5962911654"; ip="192.168.0.721"},
    @{r=24; id=10023; name="Machine 23"; serial="FB5962911655"; ip="192.168.0.841"},
    @{r=25; id=10024; name="Machine 24"; serial="FB5962911656"; ip="192.168.0.186"},
    @{r=26; id=10025; name="Machine 25"; serial="FB5962911657"; ip="192.168.0.627"},
    @{r=27; id=10026; name="Machine 26"; serial="FB5962911658"; ip="192.168.0.879"},
    @{r=28; id=10027; name="Machine 27"; serial="FB5962911659"; ip="192.168.0.628"},
    @{r=29; id=10028; name="Machine 28"; serial="FB5962911661"; ip="192.168.0.306"},
    @{r=30; id=10029; name="Machine 29"; serial="FB5962911662"; ip="192.168.0.355"}
)

foreach ($row in $newRows) {
    $r = $row.r
    $ws.Range("A$r").Value = $row.id
    $ws.Range("B$r").Value = $row.name
    $ws.Range("D$r").Value = $row.serial
    $ws.Range("E$r").Value = $row.ip
    $ws.Range("F$r").Value = 1001
    $ws.Range("G$r").Value = "eng"
    $ws.Range("H$r").Value = $true
    $ws.Range("I$r").Value = "superadmin"
    $ws.Range("J$r").Value = "now()"
    $ws.Range("K$r").Value = "now()"
}

# ---------------------------------------------------------------------------
# 2. Replace every mac_address value (rows 2-30) with the new
#    dash-separated format, top to bottom.
# ---------------------------------------------------------------------------
$macs = @(
    "8C-16-45-5A-5D-0D",
    "8C-16-45-88-E1-0D",
    "00-FF-D3-E3-9A-27",
    "8C-16-45-5A-62-41",
    "E8-6A-64-1D-75-E4",
    "8C-16-45-FA-94-B7",
    "8C-16-45-1A-0F-62",
    "E8-6A-64-1C-52-6E",
    "48-51-B7-10-35-A6",
    "8C-16-45-38-F3-F3",
    "D4-3D-7E-58-CC-45",
    "8C-16-45-5A-5D-96",
    "8C-16-45-5A-5D-8E",
    "8C-16-45-33-A5-5F",
    "3C-95-09-F9-EA-DF",
    "8C-16-45-88-E7-0B",
    "B4-69-21-5A-DB-C4",
    "E8-6A-64-1D-48-B7",
    "8C-16-45-59-69-09 ",
    "98-E7-F4-30-16-5A ",
    "38-BA-F8-53-C7-8F",
    "E8-6A-64-1C-58-C2",
    "E4-A4-71-CE-BA-93",
    "54-E1-AD-EA-30-C9",
    "8C-16-45-65-DD-40",
    "58-20-B1-D6-C3-BE",
    "8C-16-45-38-F0-25",
    "6C-88-14-AC-EF-55",
    "3C-6A-A7-C0-DF-27"
)

for ($i = 0; $i -lt $macs.Length; $i++) {
    $r = $i + 2
    $ws.Range("C$r").Value = $macs[$i]
}

# ---------------------------------------------------------------------------
# 3. Cosmetic sheet-level updates that accompanied the data edit:
#    - mac_address column widened to fit the longer dash-formatted values
#    - selection left on the row below the new data (whole-row style select)
# ---------------------------------------------------------------------------
$ws.Columns("C").ColumnWidth = 17 - 5/6

$ws.Range("A31:XFD1048576").Select() | Out-Null
